# Update points for project
# - Remove the "AutoGrading Rubric" header row from the Rubric sheet
#   (everything below it shifts up by one row).
# - Collapse the old two-section rubric (AutoGrading + Manual Grading)
#   into a single section: the old "Sub total" row becomes "Total Points"
#   (bold) with its SUM formula adjusted, and the old "Manual Grading
#   Rubric" section's rows are cleared out, leaving blank formatted rows
#   behind (matching the trailing blank rows already on the sheet).
# - Update the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rubric")

# Delete the "AutoGrading Rubric" header row (old row 2). Excel shifts
# every row below it up by one and automatically repairs formula
# references (e.g. SUM(C3:C9) -> SUM(C2:C8)).
$ws.Rows.Item(2).Delete()

# The old "Sub total" row (now row 9, holding =SUM(C2:C8)) becomes the
# sheet's single "Total Points" row: bold label, same total formula.
$ws.Range("A9").Value = "Total Points"
$ws.Range("A9").Font.Bold = $true

# The remainder of the old "Manual Grading Rubric" section (now rows
# 10-17) is no longer used - clear its leftover labels/points, leaving
# blank (but still formatted) rows like the rest of the sheet.
$ws.Range("A10:A17").ClearContents()
$ws.Range("C10:C17").ClearContents()

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("A1:C9").Select()
